$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Capture the original (un-styled) style of column D up front, before any
# writes happen, so that re-applying a style after a text-forcing write
# always restores the true original appearance (avoids picking up a
# style that this script itself just changed on a neighboring row).
$origStyleD = $ws.Range("D2").Style

# Simple text / percentage / non-numeric-looking price updates
$ws.Range("D2").Value = "39.768.18"
$ws.Range("E2").Value = "  -0.87%  "
$ws.Range("D3").Value = "2.196.20"
$ws.Range("E3").Value = "  -2.01%  "
$ws.Range("E4").Value = "  +0.00%  "
$ws.Range("E5").Value = "  -0.87%  "
$ws.Range("E6").Value = "  -1.34%  "
$ws.Range("E7").Value = "  -1.91%  "
$ws.Range("E8").Value = "  +0.12%  "
$ws.Range("E9").Value = "  -2.66%  "
$ws.Range("E10").Value = "  -3.81%  "
$ws.Range("E11").Value = "  +6.44%  "
$ws.Range("E13").Value = "  +2.42%  "
$ws.Range("E14").Value = "  -0.21%  "
$ws.Range("D15").Value = "2.536.06"
$ws.Range("E15").Value = "  -1.91%  "
$ws.Range("E16").Value = "  -3.42%  "
$ws.Range("D17").Value = "2.199.17"
$ws.Range("E17").Value = "  -1.76%  "
$ws.Range("E18").Value = "  -1.10%  "
$ws.Range("D19").Value = "39.660.58"
$ws.Range("E19").Value = "  -0.94%  "
$ws.Range("D20").Value = "0.0₃0880"
$ws.Range("E20").Value = "  -1.09%  "
$ws.Range("E21").Value = "  -1.14%  "
$ws.Range("E22").Value = "  -2.58%  "
$ws.Range("E23").Value = "  -1.27%  "
$ws.Range("E24").Value = "  +0.32%  "
$ws.Range("E25").Value = "  +0.08%  "
$ws.Range("E26").Value = "  -1.97%  "
$ws.Range("E27").Value = "  -2.97%  "
$ws.Range("E28").Value = "  +1.24%  "
$ws.Range("E29").Value = "  -3.43%  "
$ws.Range("E30").Value = "  -2.01%  "
$ws.Range("E31").Value = "  +2.66%  "
$ws.Range("E32").Value = "  -6.39%  "
$ws.Range("E33").Value = "  -0.03%  "
$ws.Range("E34").Value = "  -1.16%  "
$ws.Range("E35").Value = "  -2.87%  "
$ws.Range("E36").Value = "  -2.32%  "
$ws.Range("E37").Value = "  +0.25%  "
$ws.Range("E38").Value = "  -0.36%  "
$ws.Range("E39").Value = "  -3.16%  "
$ws.Range("B40").Value = "Celestia"
$ws.Range("C40").Value = "https://coinranking.com/coin/YQcD0lBl7+celestia-tia"
$ws.Range("E40").Value = "  -7.39%  "
$ws.Range("B41").Value = "ARBITRUM"
$ws.Range("C41").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("E41").Value = "  -2.85%  "
$ws.Range("D42").Value = "2.114.49"
$ws.Range("E42").Value = "  +2.57%  "
$ws.Range("E43").Value = "  -3.45%  "
$ws.Range("E44").Value = "  -0.80%  "
$ws.Range("E46").Value = "  -1.14%  "
$ws.Range("E47").Value = "  -5.27%  "
$ws.Range("E48").Value = "  +2.08%  "
$ws.Range("D49").Value = "2.400.32"
$ws.Range("E49").Value = "  -1.75%  "
$ws.Range("E50").Value = "  +0.30%  "
$ws.Range("E51").Value = "  -1.77%  "

# Price cells whose new text looks like a plain decimal number.
# Prefix with an apostrophe to force text, then restore the original
# (un-styled) cell style so no stray number-format style is introduced.
$ws.Range("D5").Value = "'290.67"
$ws.Range("D5").Style = $origStyleD
$ws.Range("D6").Value = "'86.06"
$ws.Range("D6").Style = $origStyleD
$ws.Range("D7").Value = "'0.507"
$ws.Range("D7").Style = $origStyleD
$ws.Range("D9").Value = "'0.464"
$ws.Range("D9").Style = $origStyleD
$ws.Range("D11").Value = "'50.00"
$ws.Range("D11").Style = $origStyleD
$ws.Range("D14").Value = "'6.41"
$ws.Range("D14").Style = $origStyleD
$ws.Range("D16").Value = "'13.67"
$ws.Range("D16").Style = $origStyleD
$ws.Range("D18").Value = "'0.726"
$ws.Range("D18").Style = $origStyleD
$ws.Range("D21").Value = "'11.14"
$ws.Range("D21").Style = $origStyleD
$ws.Range("D22").Value = "'5.70"
$ws.Range("D22").Style = $origStyleD
$ws.Range("D23").Value = "'65.08"
$ws.Range("D23").Style = $origStyleD
$ws.Range("D24").Value = "'237.23"
$ws.Range("D24").Style = $origStyleD
$ws.Range("D28").Value = "'23.25"
$ws.Range("D28").Style = $origStyleD
$ws.Range("D30").Value = "'9.15"
$ws.Range("D30").Style = $origStyleD
$ws.Range("D31").Value = "'155.76"
$ws.Range("D31").Style = $origStyleD
$ws.Range("D32").Value = "'31.21"
$ws.Range("D32").Style = $origStyleD
$ws.Range("D34").Value = "'4.91"
$ws.Range("D34").Style = $origStyleD
$ws.Range("D37").Value = "'2.84"
$ws.Range("D37").Style = $origStyleD
$ws.Range("D39").Value = "'0.0972"
$ws.Range("D39").Style = $origStyleD
$ws.Range("D40").Value = "'15.12"
$ws.Range("D40").Style = $origStyleD
$ws.Range("D41").Value = "'1.67"
$ws.Range("D41").Style = $origStyleD
$ws.Range("D43").Value = "'3.71"
$ws.Range("D43").Style = $origStyleD
$ws.Range("D44").Value = "'0.0267"
$ws.Range("D44").Style = $origStyleD
$ws.Range("D45").Value = "'2.09"
$ws.Range("D45").Style = $origStyleD
$ws.Range("D46").Value = "'9.71"
$ws.Range("D46").Style = $origStyleD
$ws.Range("D47").Value = "'17.22"
$ws.Range("D47").Style = $origStyleD
$ws.Range("D48").Value = "'2.66"
$ws.Range("D48").Style = $origStyleD
$ws.Range("D51").Value = "'88.03"
$ws.Range("D51").Style = $origStyleD
